$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.05045533333333333
$ws.Range("H2").Value = 0.151366
$ws.Range("I2").Value = 0.004442474524580737
$ws.Range("J2").Value = 0.004442474524580737
$ws.Range("M2").Value = 14.321881
$ws.Range("N2").Value = 42.965643
$ws.Range("O2").Value = 0.2949569176783066
$ws.Range("P2").Value = 0.2949569176783066
$ws.Range("Q2").Value = 0.7226152798153332
$ws.Range("R2").Value = 6.503537518338
$ws.Range("S2").Value = 0.001310338592634735
$ws.Range("T2").Value = 0.001310338592634735
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.05045533333333333
$ws.Range("H3").Value = 0.151366
$ws.Range("I3").Value = 0.004442474524580737
$ws.Range("J3").Value = 0.004442474524580737
$ws.Range("N3").Value = 81.25250700000001
$ws.Range("O3").Value = 0.557794259435499
$ws.Range("P3").Value = 0.557794259435499
$ws.Range("Q3").Value = 1.366540774951333
$ws.Range("R3").Value = 12.298866974562
$ws.Range("S3").Value = 0.002477986787499583
$ws.Range("T3").Value = 0.002477986787499583
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.05045533333333333
$ws.Range("H4").Value = 0.151366
$ws.Range("I4").Value = 0.004442474524580737
$ws.Range("J4").Value = 0.004442474524580737
$ws.Range("M4").Value = 7.149790333333333
$ws.Range("N4").Value = 21.449371
$ws.Range("O4").Value = 0.1472488228861944
$ws.Range("P4").Value = 0.1472488228861943
$ws.Range("Q4").Value = 0.3607450545317777
$ws.Range("R4").Value = 3.246705490786
$ws.Range("S4").Value = 0.0006541491444464196
$ws.Range("T4").Value = 0.0006541491444464194
$ws.Range("I5").Value = 0.7425623198471305
$ws.Range("J5").Value = 0.7425623198471305
$ws.Range("M5").Value = 14.321881
$ws.Range("N5").Value = 42.965643
$ws.Range("O5").Value = 0.2949569176783066
$ws.Range("P5").Value = 0.2949569176783066
$ws.Range("Q5").Value = 120.7855837028796
$ws.Range("R5").Value = 1087.070253325917
$ws.Range("S5").Value = 0.2190238930461624
$ws.Range("T5").Value = 0.2190238930461624
$ws.Range("I6").Value = 0.7425623198471305
$ws.Range("J6").Value = 0.7425623198471305
$ws.Range("N6").Value = 81.25250700000001
$ws.Range("O6").Value = 0.557794259435499
$ws.Range("P6").Value = 0.557794259435499
$ws.Range("S6").Value = 0.4141969992838363
$ws.Range("T6").Value = 0.4141969992838363
$ws.Range("I7").Value = 0.7425623198471305
$ws.Range("J7").Value = 0.7425623198471305
$ws.Range("M7").Value = 7.149790333333333
$ws.Range("N7").Value = 21.449371
$ws.Range("O7").Value = 0.1472488228861944
$ws.Range("P7").Value = 0.1472488228861943
$ws.Range("Q7").Value = 60.29875536354988
$ws.Range("R7").Value = 542.6887982719489
$ws.Range("S7").Value = 0.1093414275171317
$ws.Range("T7").Value = 0.1093414275171317
$ws.Range("G8").Value = 2.873389
$ws.Range("H8").Value = 8.620167
$ws.Range("I8").Value = 0.2529952056282888
$ws.Range("J8").Value = 0.2529952056282888
$ws.Range("M8").Value = 14.321881
$ws.Range("N8").Value = 42.965643
$ws.Range("O8").Value = 0.2949569176783066
$ws.Range("P8").Value = 0.2949569176783066
$ws.Range("Q8").Value = 41.152335324709
$ws.Range("R8").Value = 370.371017922381
$ws.Range("S8").Value = 0.07462268603950942
$ws.Range("T8").Value = 0.07462268603950942
$ws.Range("G9").Value = 2.873389
$ws.Range("H9").Value = 8.620167
$ws.Range("I9").Value = 0.2529952056282888
$ws.Range("J9").Value = 0.2529952056282888
$ws.Range("N9").Value = 81.25250700000001
$ws.Range("O9").Value = 0.557794259435499
$ws.Range("P9").Value = 0.557794259435499
$ws.Range("Q9").Value = 77.82335327874101
$ws.Range("R9").Value = 700.4101795086691
$ws.Range("S9").Value = 0.1411192733641631
$ws.Range("T9").Value = 0.1411192733641631
$ws.Range("G10").Value = 2.873389
$ws.Range("H10").Value = 8.620167
$ws.Range("I10").Value = 0.2529952056282888
$ws.Range("J10").Value = 0.2529952056282888
$ws.Range("M10").Value = 7.149790333333333
$ws.Range("N10").Value = 21.449371
$ws.Range("O10").Value = 0.1472488228861944
$ws.Range("P10").Value = 0.1472488228861943
$ws.Range("Q10").Value = 20.54412889610633
$ws.Range("R10").Value = 184.897160064957
$ws.Range("S10").Value = 0.03725324622461623
$ws.Range("T10").Value = 0.03725324622461622
